$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet contains a series of weekly price-report rows (9-17) that share
# the same Mercado/Región/Producto metadata and only differ by date, volume
# and price columns. A new weekly entry needs to be inserted right after
# row 9 (becoming the new row 10), pushing the existing rows 10-17 down to
# 11-18.
#
# Copy row 9 (which has the shared template: styles + all the constant
# column values) and insert the copy at row 10. This shifts rows 10-17 down
# to 11-18 automatically, carrying along their original values and styles.
$ws.Rows("9:9").Copy()
$ws.Rows("10:10").Insert()

# Now overwrite the new row 10 with its own date/volume/price data.
$ws.Range("D10").Value = 44762
$ws.Range("M10").Value = 50
$ws.Range("N10").Value = 2300
$ws.Range("O10").Value = 2300
$ws.Range("P10").Value = 2300
$ws.Range("S10").Value = 2300
